$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 246, shifting existing rows 246:299 down to 247:300.
$ws.Range("A246").EntireRow.Insert()

# Populate the newly inserted row 246 with the new record.
$ws.Range("A246").Value = 10
$ws.Range("B246").Value = "Vega Modelo de Temuco"
$ws.Range("C246").Value = "La Araucanía"
$ws.Range("D246").Value = 44889
$ws.Range("E246").Value = 9
$ws.Range("F246").Value = 100112039
$ws.Range("G246").Value = "Ciboulette"
$ws.Range("H246").Value = "Sin especificar"
$ws.Range("I246").Value = "Primera"
$ws.Range("J246").Value = 85
$ws.Range("K246").Value = 5000
$ws.Range("L246").Value = 5000
$ws.Range("M246").Value = 5000
$ws.Range("N246").Value = "$/docena de atados"
$ws.Range("O246").Value = "Provincia de Cautín"
$ws.Range("P246").Value = 1667
$ws.Range("Q246").Value = 3
$ws.Range("R246").Value = "Hortaliza"
